$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows (2025-07-19 / 2025-07-20 MLS fixtures) to append starting at row 338.
# Columns: A Fecha, B Local, C Visita, D Goles Local, E Goles Visita, F Fixture ID,
#          G Corners Local, H Corners Visita, I Amarillas Local, J Amarillas Visita,
#          K Rojas Local, L Rojas Visita, M Goles 1T Local, N Goles 1T Visita,
#          O Goles 2T Local, P Goles 2T Visita, Q Posesion Local (%), R Posesion Visita (%),
#          S Resultado
$data = @(
    @('2025-07-19', 'New York Red Bulls', 'Inter Miami', 1, 5, 1326555, 2, 4, 2, 3, 0, 0, 0, 0, 1, 5, '37%', '63%', 'V'),
    @('2025-07-19', 'Atlanta United FC', 'Charlotte', 2, 3, 1326551, 7, 2, 1, 1, 0, 0, 0, 0, 2, 3, '54%', '46%', 'V'),
    @('2025-07-19', 'New England Revolution', 'Orlando City SC', 1, 2, 1326554, 9, 6, 2, 2, 0, 0, 0, 0, 1, 2, '62%', '38%', 'V'),
    @('2025-07-19', 'CF Montreal', 'Chicago Fire', 0, 2, 1326553, 5, 3, 4, 3, 0, 0, 0, 0, 0, 2, '51%', '49%', 'V'),
    @('2025-07-20', 'Columbus Crew', 'DC United', 2, 1, 1326552, 7, 0, 1, 2, 0, 1, 0, 0, 2, 1, '64%', '36%', 'L'),
    @('2025-07-20', 'Seattle Sounders', 'San Jose Earthquakes', 3, 2, 1326560, 5, 2, 1, 3, 0, 0, 0, 0, 3, 2, '53%', '47%', 'L'),
    @('2025-07-20', 'FC Dallas', 'St. Louis City', 3, 0, 1326556, 2, 4, 2, 0, 0, 0, 0, 0, 3, 0, '34%', '66%', 'L'),
    @('2025-07-20', 'Houston Dynamo', 'Philadelphia Union', 1, 1, 1326557, 2, 4, 3, 4, 0, 1, 0, 0, 1, 1, '61%', '39%', 'E'),
    @('2025-07-20', 'Sporting Kansas City', 'New York City FC', 1, 1, 1326558, 7, 4, 2, 0, 0, 0, 0, 0, 1, 1, '39%', '61%', 'E'),
    @('2025-07-20', 'Nashville SC', 'Toronto FC', 1, 0, 1326559, 6, 2, 3, 2, 0, 0, 0, 0, 1, 0, '52%', '48%', 'L'),
    @('2025-07-20', 'Real Salt Lake', 'FC Cincinnati', 0, 1, 1326561, 3, 5, 2, 4, 0, 0, 0, 0, 0, 1, '57%', '43%', 'V'),
    @('2025-07-20', 'Los Angeles FC', 'Los Angeles Galaxy', 3, 3, 1326562, 4, 3, 3, 4, 1, 0, 0, 0, 3, 3, '35%', '65%', 'E'),
    @('2025-07-20', 'Portland Timbers', 'Minnesota United FC', 1, 1, 1326563, 6, 3, 1, 3, 0, 0, 0, 0, 1, 1, '64%', '36%', 'E'),
    @('2025-07-20', 'San Diego', 'Vancouver Whitecaps', 1, 1, 1326732, 10, 2, 1, 3, 0, 0, 0, 0, 1, 1, '62%', '38%', 'E'),
)

$startRow = 338

# Columns A (dates like "2025-07-19") and Q/R (percentages like "37%") look numeric to
# Excel's input parser, so format them as Text first to keep them as literal strings
# instead of being converted to date serials / percentage numbers.
$lastRow = $startRow + $data.Length - 1
$ws.Range("A$startRow" + ":A$lastRow").NumberFormat = "@"
$ws.Range("Q$startRow" + ":R$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
